$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7040336728096008
$ws.Range("B1").Value = 1.34195613861084
$ws.Range("C1").Value = 3.826931953430176
$ws.Range("D1").Value = 2.46489143371582
$ws.Range("E1").Value = 0.6704949140548706
